$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 139.2
$ws.Range("J2").Value = 303
$ws.Range("L2").Value = 303
$ws.Range("N2").Value = -529

$ws.Range("H29").Value = 35
$ws.Range("I29").Value = 35
$ws.Range("K29").Value = 105
$ws.Range("M29").Value = 176

$ws.Range("H43").Value = 3688.5
$ws.Range("I43").Value = 2786.8572
$ws.Range("J43").Value = 10000
$ws.Range("K43").Value = 2786.8572
$ws.Range("L43").Value = 10000
$ws.Range("M43").Value = -2717.8572
$ws.Range("N43").Value = -10138

$ws.Range("H51").Value = 9515.362999999999
$ws.Range("I51").Value = 10490
$ws.Range("J51").Value = 9298.777
$ws.Range("K51").Value = 10490
$ws.Range("L51").Value = 9298.777
$ws.Range("M51").Value = -10006
$ws.Range("N51").Value = -10266.777

$ws.Range("H58").Value = 1567.7142
$ws.Range("I58").Value = 493.5
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 1480.5
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -1330.5
$ws.Range("N58").Value = -9300

$ws.Range("H62").Value = 6673123.5
$ws.Range("J62").Value = 4197
$ws.Range("L62").Value = 4197
$ws.Range("N62").Value = -5445

$ws.Range("H65").Value = 6673123.5
$ws.Range("J65").Value = 4197
$ws.Range("L65").Value = 20985
$ws.Range("N65").Value = -27225

$ws.Range("H87").Value = 68749.5
$ws.Range("J87").Value = 68749.5
$ws.Range("L87").Value = 68749.5
$ws.Range("N87").Value = -71245.5

$ws.Range("H90").Value = 68749.5
$ws.Range("J90").Value = 68749.5
$ws.Range("L90").Value = 206248.5
$ws.Range("N90").Value = -218728.5

$ws.Range("H132").Value = 2539.5
$ws.Range("I132").Value = 2027.8636
$ws.Range("K132").Value = 6083.5908
$ws.Range("M132").Value = -3553.5908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1850.1305
$ws.Range("I2").Value = 1773.9524
$ws.Range("K2").Value = 1773.9524
$ws.Range("M2").Value = -1660.9524

$ws.Range("H32").Value = 1648.1
$ws.Range("I32").Value = 1652.6632
$ws.Range("J32").Value = 1424.5
$ws.Range("K32").Value = 1652.6632
$ws.Range("L32").Value = 1424.5
$ws.Range("M32").Value = -1365.6632
$ws.Range("N32").Value = -1998.5

$ws.Range("H45").Value = 2104.95
$ws.Range("I45").Value = 2057.8948
$ws.Range("J45").Value = 2999
$ws.Range("K45").Value = 2057.8948
$ws.Range("L45").Value = 2999
$ws.Range("M45").Value = -1680.8948
$ws.Range("N45").Value = -3753

$ws.Range("H46").Value = 3645.5715
$ws.Range("I46").Value = 4790.3335
$ws.Range("J46").Value = 2787
$ws.Range("K46").Value = 4790.3335
$ws.Range("L46").Value = 2787
$ws.Range("M46").Value = -4471.3335
$ws.Range("N46").Value = -3425

$ws.Range("H61").Value = 2819.2983
$ws.Range("I61").Value = 2091.413
$ws.Range("K61").Value = 2091.413
$ws.Range("M61").Value = -1879.413

$ws.Range("H110").Value = 3893.1428
$ws.Range("I110").Value = 4205.6665
$ws.Range("J110").Value = 2955.5715
$ws.Range("K110").Value = 4205.6665
$ws.Range("L110").Value = 2955.5715
$ws.Range("M110").Value = -2160.6665
$ws.Range("N110").Value = -7045.5715

$ws.Range("H116").Value = 1850.1305
$ws.Range("I116").Value = 1773.9524
$ws.Range("K116").Value = 1773.9524
$ws.Range("M116").Value = 520.0476000000001

$ws.Range("H132").Value = 2974.077
$ws.Range("I132").Value = 2974.077
$ws.Range("K132").Value = 8922.231
$ws.Range("M132").Value = -6392.231

$ws.Range("H136").Value = 2819.2983
$ws.Range("I136").Value = 2091.413
$ws.Range("K136").Value = 6274.239
$ws.Range("M136").Value = -3724.239

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1850.1305
$ws.Range("I3").Value = 1773.9524
$ws.Range("K3").Value = 1773.9524
$ws.Range("M3").Value = -1659.9524

$ws.Range("H105").Value = 1821.5333
$ws.Range("I105").Value = 1758.625
$ws.Range("K105").Value = 1758.625
$ws.Range("M105").Value = -11.625

$ws.Range("H140").Value = 93876.664
$ws.Range("J140").Value = 93876.664
$ws.Range("L140").Value = 93876.664
$ws.Range("N140").Value = -104236.664

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 196832.67
$ws.Range("I6").Value = 294999.5
$ws.Range("J6").Value = 499
$ws.Range("K6").Value = 294999.5
$ws.Range("L6").Value = 499
$ws.Range("M6").Value = -294886.5
$ws.Range("N6").Value = -725

$ws.Range("H16").Value = 1719.25
$ws.Range("J16").Value = 10013
$ws.Range("L16").Value = 10013
$ws.Range("N16").Value = -10587

$ws.Range("H19").Value = 26754496
$ws.Range("I19").Value = 33442994
$ws.Range("K19").Value = 33442994
$ws.Range("M19").Value = -33442824

$ws.Range("H24").Value = 26754496
$ws.Range("I24").Value = 33442994
$ws.Range("K24").Value = 33442994
$ws.Range("M24").Value = -33442824

$ws.Range("H37").Value = 932.3333
$ws.Range("I37").Value = 800
$ws.Range("J37").Value = 998.5
$ws.Range("K37").Value = 800
$ws.Range("L37").Value = 998.5
$ws.Range("M37").Value = -693
$ws.Range("N37").Value = -1212.5

$ws.Range("H68").Value = 36538.152
$ws.Range("J68").Value = 36538.152
$ws.Range("L68").Value = 36538.152
$ws.Range("N68").Value = -38036.152

$ws.Range("H71").Value = 36538.152
$ws.Range("J71").Value = 36538.152
$ws.Range("L71").Value = 109614.456
$ws.Range("N71").Value = -117102.456

$ws.Range("H94").Value = 1323
$ws.Range("J94").Value = 1531.6
$ws.Range("L94").Value = 1531.6
$ws.Range("N94").Value = -2433.6

$ws.Range("H113").Value = 1719.25
$ws.Range("J113").Value = 10013
$ws.Range("L113").Value = 10013
$ws.Range("N113").Value = -14353

$ws.Range("H132").Value = 4594.22
$ws.Range("I132").Value = 4487.522
$ws.Range("K132").Value = 13462.566
$ws.Range("M132").Value = -10932.566

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 132.375
$ws.Range("J11").Value = 359.5
$ws.Range("L11").Value = 1078.5
$ws.Range("N11").Value = -1358.5

$ws.Range("H40").Value = 81.43478399999999
$ws.Range("I40").Value = 88.166664
$ws.Range("J40").Value = 57.2
$ws.Range("K40").Value = 352.666656
$ws.Range("L40").Value = 228.8
$ws.Range("M40").Value = -283.666656
$ws.Range("N40").Value = -366.8

$ws.Range("H80").Value = 3064.6667
$ws.Range("J80").Value = 3166.6667
$ws.Range("L80").Value = 9500.000100000001
$ws.Range("N80").Value = -11372.0001

$ws.Range("H83").Value = 3064.6667
$ws.Range("J83").Value = 3166.6667
$ws.Range("L83").Value = 28500.0003
$ws.Range("N83").Value = -37860.0003

$ws.Range("H140").Value = 10735.261
$ws.Range("J140").Value = 5288.8335
$ws.Range("L140").Value = 15866.5005
$ws.Range("N140").Value = -26226.5005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4169015
$ws.Range("I3").Value = 3384.375
$ws.Range("J3").Value = 12500276
$ws.Range("K3").Value = 3384.375
$ws.Range("L3").Value = 12500276
$ws.Range("M3").Value = -3268.375
$ws.Range("N3").Value = -12500508

$ws.Range("H10").Value = 15005000
$ws.Range("I10").Value = 15005000
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 15005000
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -15004831
$ws.Range("N10").ClearContents()

$ws.Range("H35").Value = 20000
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20596

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 1000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H16").Value = 3156.5715
$ws.Range("I16").Value = 2332
$ws.Range("J16").Value = 3775
$ws.Range("K16").Value = 2332
$ws.Range("L16").Value = 3775
$ws.Range("M16").Value = -2162
$ws.Range("N16").Value = -4115

$ws.Range("H103").Value = 19999
$ws.Range("J103").Value = 19999
$ws.Range("L103").Value = 19999
$ws.Range("N103").Value = -22343

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 35622.5
$ws.Range("I2").Value = 39166
$ws.Range("K2").Value = 39166
$ws.Range("M2").Value = -39054

$ws.Range("H100").Value = 1872.25
$ws.Range("J100").Value = 4999
$ws.Range("L100").Value = 9998
$ws.Range("N100").Value = -11080

$ws.Range("H132").Value = 1840.0454
$ws.Range("I132").Value = 1244.8718
$ws.Range("K132").Value = 3734.6154
$ws.Range("M132").Value = -1204.6154

$ws.Range("H136").Value = 324817.8
$ws.Range("I136").Value = 335611.78
$ws.Range("K136").Value = 1006835.34
$ws.Range("M136").Value = -1004285.34
